# Probability higher revision in costs and utilities
#
# Inserts a new "other_rates" worksheet between "utilities" and
# "uk_lifetables" holding the log-rate of higher (3rd+) revision mean/SE
# used for the probability-of-higher-revision calculation.

$wb = $excel.ActiveWorkbook

# Insert the new sheet directly after "utilities" so it lands in the same
# slot as the 4th tab (before "uk_lifetables").
$afterSheet = $wb.Worksheets.Item("utilities")
$ws = $wb.Worksheets.Add($null, $afterSheet)
$ws.Name = "other_rates"

# Header row
$ws.Range("A1").Value = "parameter"
$ws.Range("B1").Value = "value"
$ws.Range("C1").Value = "comment"

# Data rows
$ws.Range("A2").Value = "lograte_higher_revision_mean"
$ws.Range("A3").Value = "lograte_higher_revision_se"

$ws.Range("B2").Value = -3.073387
$ws.Range("B3").Value = 0.0001499125

$comment = "Inverse variance meta-analysis of log rates 3rd,.., 8th revision from NJR estimates"
$ws.Range("C2").Value = $comment
$ws.Range("C3").Value = $comment

# Column widths roughly matching the authored layout
$ws.Columns.Item(1).ColumnWidth = 24.333333333333332
$ws.Columns.Item(2).ColumnWidth = 17

$ws.PageSetup.Orientation = 1

$null = $ws.Range("C4").Select()
